$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (14) on the "Repayment Schedule"
# sheet. This shifts the old N/O/P columns (header "Late" / blank / header
# "Outstanding" and their per-row values) one column to the right, into
# O/P/Q, leaving the new N column blank - matching the source workbook's
# move from "mifos" to "finflux" column layout.
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab and move the selection to
# R6 (this also clears the previous tabSelected flag on "Transactions" and
# updates the workbook's activeTab).
$ws.Activate()
$ws.Range("R6").Select()
